# Q3 Update - 2025
# Updates a set of numeric-text cells (stored as shared strings, not real
# numbers) in the "fromCSV" sheet. Values are copied in from cells that
# already contain the desired text (as plain numeric-looking strings) so
# that the destination keeps its original cell style and the shared string
# table keeps reusing existing entries, exactly like the source data set.
# A couple of brand-new numeric strings are staged through a scratch cell
# (formatted as Text) and pasted in as values-only so the destination style
# is left untouched while the new text is still stored as a string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- direct "value already exists elsewhere as text" swaps -----------------
# map: destination cell -> source cell already holding the desired text
$copies = @{
    "O488" = "D12";   # asylum_seekers 9   -> 11
    "O490" = "D25";   # asylum_seekers 19  -> 24
    "P490" = "D6";    # returned_refugees 0 -> 5
    "N491" = "D56";   # refugees 51  -> 55
    "N493" = "D20";   # refugees 18  -> 19
    "N494" = "D23";   # refugees 20  -> 22
    "O496" = "F28";   # asylum_seekers 420 -> 41
    "N497" = "D323";  # refugees 0   -> 322
    "O500" = "D473";  # asylum_seekers 471 -> 472
    "N501" = "D38";   # refugees 36  -> 37
}

foreach ($dest in $copies.Keys) {
    $src = $copies[$dest]
    $ws.Range($src).Copy()
    $ws.Range($dest).PasteSpecial(-4163)   # xlPasteValues
}

# --- brand new numeric-text values -----------------------------------------
# staged through a scratch cell (forced to Text) so they remain strings
# rather than being auto-converted to numbers, then pasted as values-only
# so the target cell keeps its own existing style.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"

$newValues = @{
    "N496" = "39265";  # refugees 38149 -> 39265
    "P496" = "541";    # returned_refugees 395 -> 541
    "O497" = "930";    # asylum_seekers 1037 -> 930
}

foreach ($dest in $newValues.Keys) {
    $scratch.Value = $newValues[$dest]
    $scratch.Copy()
    $ws.Range($dest).PasteSpecial(-4163)   # xlPasteValues
}

$scratch.Clear()

$excel.CutCopyMode = 0

# --- rename the sheet's short-url value (shared by every row) --------------
$ws.Cells.Replace("V7tm0w", "6m5DqR")
